# Apply "update resource by parse tool" edit to the 玩家-商城配置 (player-shop)
# worksheet: fill in the two newly-parsed shop rows (effect / foot) that sit
# below the existing "hero" row, flip their purchase-limit flags on, and
# move the sheet's active selection down to the newly added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 ("hero" entry) -------------------------------------------------
# Text values are unchanged; re-assert them for clarity/safety. The real
# change on this row is BuyLimit (L5) flipping from 0 to 1.
$ws.Range("D5").Value = "hero"
$ws.Range("E5").Value = '[{"hero":{"id":"3"}}]'
$ws.Range("G5").Value = 'money=[{"money":"90"}]'
$ws.Range("L5").Value = 1

# --- Row 6 ("effect" entry, newly populated) ------------------------------
$ws.Range("B6").Value = 11001
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "effect"
$ws.Range("E6").Value = '[{"effect":{"id":"3"}}]'
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 'money=[{"money":"90"}]'
$ws.Range("L6").Value = 1

# --- Row 7 ("foot" entry, newly populated) --------------------------------
$ws.Range("B7").Value = 12001
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "foot"
$ws.Range("E7").Value = '[{"foot":{"id":"3"}}]'
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 'money=[{"money":"90"}]'
$ws.Range("L7").Value = 1

# --- View state: move the selection to the last edited cell and scroll the
# sheet so column G is the left-most visible column, matching the author's
# on-screen position after the edit.
$ws.Range("L7").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1
